# Update Excel sheet with new scrape data (2025-11-27 03:20:32 UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (columns C, D, F, H change) ---
# Target OOXML widths: C=47, D=70, F=17, H=53
# COM ColumnWidth differs from the stored OOXML width by a fixed offset (~0.83)
$ws.Columns.Item(3).ColumnWidth = 46.17
$ws.Columns.Item(4).ColumnWidth = 69.17
$ws.Columns.Item(6).ColumnWidth = 16.17
$ws.Columns.Item(8).ColumnWidth = 52.17

# --- Update data rows (header row 1 unchanged) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1330126"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330126"
$ws.Range("C2").Value = "Social Media Marketing"
$ws.Range("D2").Value = "Oporto, Portugal"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "4 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "StayWell"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1330113"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1330113"
$ws.Range("C3").Value = "[Accelerate Serbia] Interior Designer Intern"
$ws.Range("D3").Value = "Subotica, Serbia"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Studio White"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1330095"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1330095"
$ws.Range("C4").Value = "International Relations Development Intern"
$ws.Range("D4").Value = "Phagwara, Punjab, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Lovely Professional University"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "1329929"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1329929"
$ws.Range("C5").Value = "MARKETING & SALES"
$ws.Range("D5").Value = "Denizli, Kumkısık, Denizli, Türkiye"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "1 applicant"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "DOĞAN İNŞ.MLZM.HAFR.NAK.İŞ MAK.SAN.VE TİC.LTD.ŞTİ."

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "1329633"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1329633"
$ws.Range("C6").Value = "Sales Coordinator"
$ws.Range("D6").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "6 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Sodexo Mexico"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "1328650"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1328650"
$ws.Range("C7").Value = "Marketing Intern"
$ws.Range("D7").Value = "Novi Sad, Serbia"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "48 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "DataDrill"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "1328614"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1328614"
$ws.Range("C8").Value = "Field Service Engineer"
$ws.Range("D8").Value = "Madrid, Spain"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "117 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "Mitsubishi Power Europe Sucursal en España"

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1328557"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1328557"
$ws.Range("C9").Value = "Arduino Developer"
$ws.Range("D9").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "0 applicants"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "Techno square"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "1328206"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1328206"
$ws.Range("C10").Value = "Power BI Specialist"
$ws.Range("D10").Value = "Frankfurt am Main, Deutschland"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "143 applicants"
$ws.Range("G10").Value = "3 - 6 Months"
$ws.Range("H10").Value = "Greyfood GmbH"

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "1323714"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1323714"
$ws.Range("C11").Value = "Digital Marketing"
$ws.Range("D11").Value = "Oporto, Portugal"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "72 applicants"
$ws.Range("G11").Value = "9 - 12 Weeks"
$ws.Range("H11").Value = "Obras Descomplicadas Lda"

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1322463"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1322463"
$ws.Range("C12").Value = "E-commerce Sales Manager"
$ws.Range("D12").Value = "Giza, El Omraniya, Giza Governorate, Egypt"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "5 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "EG scout shop"

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "1312564"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1312564"
$ws.Range("C13").Value = "Sales Intern"
$ws.Range("D13").Value = "Delhi, India"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "2 applicants"
$ws.Range("G13").Value = "6 - 18 Months"
$ws.Range("H13").Value = "NRM International"

